$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B4").Value = 1056000000.0
$ws.Range("C4").Value = 1030000000.0
$ws.Range("D4").Value = 1064000000.0
$ws.Range("E4").Value = 1228000000.0
$ws.Range("F4").Value = 1227000000.0

$ws.Range("B13").Value = 1033000000.0
$ws.Range("C13").Value = 991000000.0
$ws.Range("D13").Value = 697000000.0
$ws.Range("E13").Value = 729000000.0
$ws.Range("F13").Value = 895000000.0

$ws.Range("B20").Value = 85000000.0
$ws.Range("C20").Value = 85000000.0
$ws.Range("D20").Value = 136000000.0
$ws.Range("E20").Value = 136000000.0
$ws.Range("F20").Value = 199000000.0

$ws.Range("G34").Value = 6384000000.0
$ws.Range("G35").Value = 7429000000.0
